# Apply cryptos list update (prices / volume changes / row re-ordering)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '50.997.93'
Set-TextCell 2 5 '  -0.15%  '

Set-TextCell 3 4 '2.955.50'
Set-TextCell 3 5 '  +0.53%  '

Set-TextCell 4 5 '  +0.00%  '

Set-TextCell 5 4 '378.55'
Set-TextCell 5 5 '  -0.27%  '

Set-TextCell 6 4 '101.90'
Set-TextCell 6 5 '  -0.50%  '

Set-TextCell 7 5 '  +0.60%  '

Set-TextCell 8 5 '  +0.00%  '

Set-TextCell 9 4 '0.584'
Set-TextCell 9 5 '  -0.47%  '

Set-TextCell 10 4 '36.28'
Set-TextCell 10 5 '  -0.77%  '

Set-TextCell 11 5 '  -0.38%  '

Set-TextCell 12 5 '  +1.07%  '

Set-TextCell 13 2 'Uniswap'
Set-TextCell 13 3 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell 13 4 '12.49'
Set-TextCell 13 5 '  +74.92%  '

Set-TextCell 14 2 'Chainlink'
Set-TextCell 14 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 14 4 '18.41'
Set-TextCell 14 5 '  +2.16%  '

Set-TextCell 15 2 'Polkadot'
Set-TextCell 15 3 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 15 4 '7.79'
Set-TextCell 15 5 '  +5.42%  '

Set-TextCell 16 2 'WrappedliquidstakedEther2.0'
Set-TextCell 16 3 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell 16 4 '3.415.63'
Set-TextCell 16 5 '  +0.23%  '

Set-TextCell 17 4 '2.959.85'
Set-TextCell 17 5 '  +2.01%  '

Set-TextCell 18 5 '  +2.68%  '

Set-TextCell 19 4 '50.973.84'
Set-TextCell 19 5 '  -0.05%  '

Set-TextCell 20 4 '3.10'
Set-TextCell 20 5 '  -3.30%  '

Set-TextCell 21 4 '12.40'
Set-TextCell 21 5 '  -1.18%  '

Set-TextCell 22 4 '0.0₃0952'
Set-TextCell 22 5 '  -0.04%  '

Set-TextCell 23 4 '69.49'
Set-TextCell 23 5 '  +1.40%  '

Set-TextCell 24 4 '266.30'
Set-TextCell 24 5 '  +1.68%  '

Set-TextCell 25 4 '3.25'
Set-TextCell 25 5 '  +11.60%  '

Set-TextCell 26 4 '8.09'
Set-TextCell 26 5 '  -2.49%  '

Set-TextCell 27 5 '  +0.04%  '

Set-TextCell 28 4 '7.02'
Set-TextCell 28 5 '  -8.49%  '

Set-TextCell 29 2 'Kaspa'
Set-TextCell 29 3 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell 29 4 '0.164'
Set-TextCell 29 5 '  -2.77%  '

Set-TextCell 30 2 'EthereumClassic'
Set-TextCell 30 3 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 30 4 '25.67'
Set-TextCell 30 5 '  +0.11%  '

Set-TextCell 31 5 '  -4.55%  '

Set-TextCell 32 4 '10.22'
Set-TextCell 32 5 '  +4.15%  '

Set-TextCell 33 4 '50.49'

Set-TextCell 34 2 'InjectiveProtocol'
Set-TextCell 34 3 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell 34 4 '33.87'
Set-TextCell 34 5 '  -1.13%  '

Set-TextCell 35 2 'Toncoin'
Set-TextCell 35 3 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 35 4 '2.05'
Set-TextCell 35 5 '  +0.20%  '

Set-TextCell 36 5 '  -5.08%  '

Set-TextCell 38 4 '3.13'
Set-TextCell 38 5 '  +4.96%  '

Set-TextCell 39 5 '  +1.00%  '

Set-TextCell 40 4 '16.57'
Set-TextCell 40 5 '  -1.42%  '

Set-TextCell 41 5 '  +2.15%  '

Set-TextCell 42 4 '2.51'
Set-TextCell 42 5 '  -2.96%  '

Set-TextCell 43 4 '118.29'
Set-TextCell 43 5 '  -2.05%  '

Set-TextCell 44 4 '3.55'
Set-TextCell 44 5 '  +9.81%  '

Set-TextCell 45 4 '21.40'
Set-TextCell 45 5 '  +0.56%  '

Set-TextCell 46 5 '  -1.65%  '

Set-TextCell 47 5 '  -3.27%  '

Set-TextCell 48 4 '2.007.93'
Set-TextCell 48 5 '  +0.01%  '

Set-TextCell 49 5 '  -5.33%  '

Set-TextCell 50 4 '0.0319'
Set-TextCell 50 5 '  -8.05%  '

Set-TextCell 51 4 '5.31'
Set-TextCell 51 5 '  +4.99%  '
